$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.08%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.25'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.29%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.184'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.48%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05943'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.35%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.704'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.49%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8685'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.66%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.002'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '8.69%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1418'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.61%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.03552'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4.77%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07187'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.46%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03148'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.63%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09250'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.29%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001543'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.17%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005998'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.70%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.490'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.31%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.260'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.80%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.227'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.25%'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'One'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.01060'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1,656.96%'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3147'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.63%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1306'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.56%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.561'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.14%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04283'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.85%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1412'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.43%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001223'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.51%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004518'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-10.48%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.11%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-22.97%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03838'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.05%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006579'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '15.60%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.97%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002200'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.11%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01049'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '7.65%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005486'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '4.01%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.11%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '28.49%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002178'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.18%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.11%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.11%'
